$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "r775"
$ws.Range("B18").Value = "richard"
$ws.Range("C18").Value = "pss doesnt mention which weight to use"
$ws.Range("D18").Value = "2025-10-01 16:14:36"
